# Databases/中文_words.xlsx : "updated DBs and added issue #22"
# Adds three new vocabulary entries to the bottom of the word list:
#   保留  -> 1.) v. preserve, hold back, retain, reserve      (2020-12-08)
#   授予  -> 1.) vt. award, confer, grand, endow              (2020-12-08)
#   诺贝尔 -> 1.) Nobel Prize                                  (2020-12-08)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last populated row currently is 88 (header is row 1), so new entries land
# on rows 89-91.
$startRow = $ws.UsedRange.Rows.Count + 1

$newEntries = @(
    @("保留",   "1.) v. preserve, hold back, retain, reserve"),
    @("授予",   "1.) vt. award, confer, grand, endow"),
    @("诺贝尔", "1.) Nobel Prize")
)

$entryDate = "2020-12-08"

# Format the new date cells as Text first so Excel stores the literal
# "2020-12-08" string (matching every other DATE cell in the sheet) instead
# of silently coercing it into a serial date number.
$dateRange = $ws.Range($ws.Cells.Item($startRow, 3), $ws.Cells.Item($startRow + $newEntries.Count - 1, 3))
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $newEntries.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newEntries[$i][0]
    $ws.Cells.Item($r, 2).Value = $newEntries[$i][1]
    $ws.Cells.Item($r, 3).Value = $entryDate
}
